$wb = $excel.ActiveWorkbook

# --- New GUIDs / content-hashes generated for this handoff report run ---
# (old a3656b9f-10d6-436e-b0ae-a9c49efaf5d9 -> new 58df128e-616d-4a68-9c4f-8d9933a956a4)
# (old ab375c47-ec3c-4984-bbf7-ebe65db1de3b -> new c231cde8-c7a8-46fc-831a-6c0baefb0bd2)
$newId1 = "58df128e-616d-4a68-9c4f-8d9933a956a4"
$newId2 = "c231cde8-c7a8-46fc-831a-6c0baefb0bd2"

# (old 0f3c5468a16b639d57cb3865365e6b815acdf216 -> new b3d4616001cec1869aa6d0af1e42184e88870e59)
# (old db4334fdc565deea208664826e9c018d2e8d045b -> new 2ec9f056696a320fe4b44c07e3caf8e61108e3a5)
$newHash1 = "b3d4616001cec1869aa6d0af1e42184e88870e59"
$newHash2 = "2ec9f056696a320fe4b44c07e3caf8e61108e3a5"

# =========================================================
# Sheet "Overview"
# =========================================================
$ws = $wb.Worksheets.Item("Overview")

# Update filename + path cells (A2/A3, B2/B3) and fix hyperlinks
$linkTarget1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1b40866264f751e896509290e256ccc66a08c6d/e2e/$newId1.md"
$linkTarget2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1b40866264f751e896509290e256ccc66a08c6d/e2e/$newId2.md"

$ws.Range("A2").Value2 = "$newId1.md"
$ws.Range("A3").Value2 = "$newId2.md"

$hls = $ws.Hyperlinks
$hls.Delete()
$hls.Add($ws.Range("B2"), $linkTarget1, "", "", "e2e\$newId1.md")
$hls.Add($ws.Range("B3"), $linkTarget2, "", "", "e2e\$newId2.md")

$ws.Range("G2").Value2 = "2016-11-29 02:26:57"
$ws.Range("G3").Value2 = "2016-11-29 02:26:57"

# =========================================================
# Sheet "zh-cn"
# =========================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("G2").Value2 = "$newId1.$newHash1.zh-cn.xlf"
$ws.Range("G3").Value2 = "$newId2.$newHash2.zh-cn.xlf"
$ws.Range("H2").Value2 = "2016-11-29 02:26:44"
$ws.Range("H3").Value2 = "2016-11-29 02:26:44"

$hls = $ws.Hyperlinks
$hls.Delete()
$hls.Add($ws.Range("A2"), $linkTarget1, "", "", "$newId1.md")
$hls.Add($ws.Range("A3"), $linkTarget2, "", "", "$newId2.md")

# =========================================================
# Sheet "de-de"
# =========================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("G2").Value2 = "$newId1.$newHash1.de-de.xlf"
$ws.Range("G3").Value2 = "$newId2.$newHash2.de-de.xlf"
$ws.Range("H2").Value2 = "2016-11-29 02:26:57"
$ws.Range("H3").Value2 = "2016-11-29 02:26:57"

$hls = $ws.Hyperlinks
$hls.Delete()
$hls.Add($ws.Range("A2"), $linkTarget1, "", "", "$newId1.md")
$hls.Add($ws.Range("A3"), $linkTarget2, "", "", "$newId2.md")

Write-Host "All updates applied"
